# Update dx uid in 24Tto25TMap to 2024 mer targets:
# column B ("dx") holds the DataElementGroup UID used for the 2023 cop_year
# rows; replace the old UID with the new 2024 mer-targets UID for every
# data row (row 1 is the header "dx").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Map")

$oldUid = "DE_GROUP-OuKFZzVk6gr"
$newUid = "DE_GROUP-TXAVaM4oYMd"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq $oldUid) {
        $cell.Value = $newUid
    }
}

# Restore the cursor/selection position recorded in the sheet view.
$ws.Range("G12").Select()
